# Daily attendance processing - 2025-11-13 18:29:54
# Swap the order of "email, System" -> "System, email" in the
# "Recorded By" column (G) for the affected attendance rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3,6,7,10,12,13,14,15,18,19,20,21,22,24,26,29,32,33,36,38,39,40,41,44,45,46,47,48,50,52,55,58,59,62,64,65,66,67,70,71,72,73,74,76,78,83,84,85,86,90,92,99,101,109,110,111,112,116,118,125,127,135,136,137,138,142,144,151,153)

foreach ($r in $rows) {
    $cell = $ws.Range("G$r")
    $current = [string]$cell.Text
    $parts = $current -split ',\s*', 2
    if ($parts.Count -eq 2) {
        $newValue = $parts[1].Trim() + ", " + $parts[0].Trim()
        $cell.Value = $newValue
    }
}
